$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each changed row.
# Price values are written as Text (NumberFormat "@") so that strings like
# "1.00" or "27.035.20" are preserved exactly instead of being parsed as numbers,
# then the style is reset to Normal so no extra formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.035.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.659.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.60%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.74%  '

$ws.Range("E6").Value = '  +1.86%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  +2.57%  '

$ws.Range("E9").Value = '  +1.78%  '

$ws.Range("E10").Value = '  +4.60%  '

$ws.Range("E11").Value = '  +4.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.892.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.57%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.654.88'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.43%  '

$ws.Range("E14").Value = '  +1.90%  '

$ws.Range("E15").Value = '  +2.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.052.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '237.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0737'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.75'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.43%  '

$ws.Range("E23").Value = '  +2.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.24%  '

$ws.Range("E26").Value = '  +1.89%  '

$ws.Range("E27").Value = '  +0.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("E30").Value = '  +0.12%  '

$ws.Range("E31").Value = '  +1.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.549.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.56%  '

$ws.Range("E33").Value = '  +1.80%  '

$ws.Range("E34").Value = '  +4.19%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.63'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.63%  '

$ws.Range("E36").Value = '  -0.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.576'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.901'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.67%  '

$ws.Range("E39").Value = '  +2.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '66.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.968'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.69%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.802.93'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.775'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.57%  '

$ws.Range("E48").Value = '  +2.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.16%  '

$ws.Range("E50").Value = '  +4.25%  '

$ws.Range("E51").Value = '  +1.02%  '
